$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new data rows right before the former row 470. Excel's
# Insert() pushes the existing rows 470:576 down to 473:579, which is
# exactly the shift seen across the whole tail of the sheet in the diff.
$ws.Rows("470:472").Insert()

# Row 470 - new weekly price point for "Provincia de Cautín"
$ws.Cells.Item(470, 1).Value = 10
$ws.Cells.Item(470, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(470, 3).Value = "La Araucanía"
$ws.Cells.Item(470, 4).Value = 45211
$ws.Cells.Item(470, 5).Value = 9
$ws.Cells.Item(470, 6).Value = 100112009
$ws.Cells.Item(470, 7).Value = "Acelga"
$ws.Cells.Item(470, 8).Value = "Sin especificar"
$ws.Cells.Item(470, 9).Value = "Primera"
$ws.Cells.Item(470, 10).Value = 90
$ws.Cells.Item(470, 11).Value = 8000
$ws.Cells.Item(470, 12).Value = 8000
$ws.Cells.Item(470, 13).Value = 8000
$ws.Cells.Item(470, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(470, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(470, 16).Value = 667
$ws.Cells.Item(470, 17).Value = 12
$ws.Cells.Item(470, 18).Value = "Hortaliza"

# Row 471 - new weekly price point for "Región Metropolitana"
$ws.Cells.Item(471, 1).Value = 10
$ws.Cells.Item(471, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(471, 3).Value = "La Araucanía"
$ws.Cells.Item(471, 4).Value = 45211
$ws.Cells.Item(471, 5).Value = 9
$ws.Cells.Item(471, 6).Value = 100112009
$ws.Cells.Item(471, 7).Value = "Acelga"
$ws.Cells.Item(471, 8).Value = "Sin especificar"
$ws.Cells.Item(471, 9).Value = "Primera"
$ws.Cells.Item(471, 10).Value = 100
$ws.Cells.Item(471, 11).Value = 7000
$ws.Cells.Item(471, 12).Value = 7000
$ws.Cells.Item(471, 13).Value = 7000
$ws.Cells.Item(471, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(471, 15).Value = "Región Metropolitana"
$ws.Cells.Item(471, 16).Value = 583
$ws.Cells.Item(471, 17).Value = 12
$ws.Cells.Item(471, 18).Value = "Hortaliza"

# Row 472 - new weekly price point for "Región del Maule"
$ws.Cells.Item(472, 1).Value = 10
$ws.Cells.Item(472, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(472, 3).Value = "La Araucanía"
$ws.Cells.Item(472, 4).Value = 45211
$ws.Cells.Item(472, 5).Value = 9
$ws.Cells.Item(472, 6).Value = 100112009
$ws.Cells.Item(472, 7).Value = "Acelga"
$ws.Cells.Item(472, 8).Value = "Sin especificar"
$ws.Cells.Item(472, 9).Value = "Primera"
$ws.Cells.Item(472, 10).Value = 80
$ws.Cells.Item(472, 11).Value = 7000
$ws.Cells.Item(472, 12).Value = 7000
$ws.Cells.Item(472, 13).Value = 7000
$ws.Cells.Item(472, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(472, 15).Value = "Región del Maule"
$ws.Cells.Item(472, 16).Value = 583
$ws.Cells.Item(472, 17).Value = 12
$ws.Cells.Item(472, 18).Value = "Hortaliza"
